# The commit swaps the embedded theme XML parts: ppt/theme/theme1.xml
# (the deck's "Integral" theme, used by the Slide Master / the deck
# itself) and ppt/theme/theme2.xml (the default "Office Theme", used by
# the Notes Master) exchange their contents. The font scheme and the
# format scheme (fills/lines/effects) are already byte-identical between
# the two theme parts, so the only semantic difference is the 12-colour
# theme colour scheme. Re-colour the deck's theme from the "Integral"
# palette to the "Office" palette so the Slide Master (and therefore the
# whole deck) ends up using the stock "Office Theme" colours.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Office Theme colour scheme, in MsoThemeColorSchemeIndex order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink.
$officeColors = @{
    1  = 0          # dk1      000000
    2  = 16777215   # lt1      FFFFFF
    3  = 6968388     # dk2      44546A
    4  = 15132391    # lt2      E7E6E6
    5  = 13998939    # accent1  5B9BD5
    6  = 3243501      # accent2  ED7D31
    7  = 10855845    # accent3  A5A5A5
    8  = 49407        # accent4  FFC000
    9  = 12874308    # accent5  4472C4
    10 = 4697456      # accent6  70AD47
    11 = 12673797    # hlink    0563C1
    12 = 7491477      # folHlink 954F72
}

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i]
}
